$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra epoch columns (F:I) entirely so the used range shrinks back to B:E
$ws.Range("F1:I2").Clear()

# Update the accuracy row with the new measured value
$ws.Range("B2").Value = 79.55729141831398
$ws.Range("C2").Value = 79.55729141831398
$ws.Range("D2").Value = 79.55729141831398
$ws.Range("E2").Value = 79.55729141831398
